# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Matches the diff: new columns AD:AF, header row styled like the rest of
# row 1 (bold, centered, bordered), data rows 2-66 filled with the
# season's W/L/T record (61-101-0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
# Clone the existing header formatting (style index used by A1:AC1) onto
# the three new header cells, then set their labels.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-66) -------------------------------------------------
# Every player row gets the team's season record.
$wins = 61
$losses = 101
$ties = 0

for ($r = 2; $r -le 66; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # AD
    $ws.Cells.Item($r, 31).Value = $losses  # AE
    $ws.Cells.Item($r, 32).Value = $ties    # AF
}
